# Add 2022-Q4 data
#
# 1. Update the "总计" (totals) summary sheet: add a new "2022-Q4" row at the
#    top of the data (row 2), shifting the existing quarters down by one row
#    and appending the oldest quarter (2020-Q4) as the new last row.
# 2. Insert a brand-new "2022-Q4" worksheet (right after "总计") holding the
#    per-fund holdings table for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" sheet — rewrite the 4-column table (日期 / 持有数量(只) / 持有市值(亿元))
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Give the brand-new row 9 the same look (border/bold/center/top, like column A
# of the other rows) by copying formatting down from row 8 before writing
# values into it.
$totalSheet.Range("A8:D8").Copy()
$totalSheet.Range("A9").PasteSpecial(-4122)

$totalRows = @(
  @(0, "2022-Q4", 2,  0.02),
  @(1, "2022-Q2", 2,  0.01),
  @(2, "2022-Q1", 3,  0.09),
  @(3, "2021-Q4", 13, 1.72),
  @(4, "2021-Q3", 8,  1.9),
  @(5, "2021-Q2", 15, 3.44),
  @(6, "2021-Q1", 22, 8.57),
  @(7, "2020-Q4", 11, 4.97)
)

for ($i = 0; $i -lt $totalRows.Length; $i++) {
  $r = $i + 2
  $row = $totalRows[$i]
  $totalSheet.Cells.Item($r, 1).Value = $row[0]
  $totalSheet.Cells.Item($r, 2).Value = $row[1]
  $totalSheet.Cells.Item($r, 3).Value = $row[2]
  $totalSheet.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------------------
# 2. New "2022-Q4" sheet — per-fund holdings table, inserted right after 总计
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
  $newSheet.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(3, 1).Value = 1

# columns B..G are stored as text (fund code keeps its leading zero, the
# numeric-looking figures keep their original string formatting)
$newSheet.Range("B2:G3").NumberFormat = "@"

$fundRows = @(
  @("014831", "兴银中证1000指数增强A", "1.40", "82.60", "1.02", "0.0143", 6),
  @("014832", "兴银中证1000指数增强C", "1.01", "82.60", "1.02", "0.0103", 6)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
  $r = $i + 2
  $row = $fundRows[$i]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $newSheet.Cells.Item($r, $c + 2).Value = $row[$c]
  }
}

# Header row + index column formatting (bold, thin border, centered/top) to
# match the look of every other quarter sheet in the workbook.
$hdrRange = $newSheet.Range("B1:H1")
$hdrRange.Font.Bold = $true
$hdrRange.HorizontalAlignment = -4108  # xlCenter
$hdrRange.VerticalAlignment = -4160    # xlTop
$hdrRange.Borders.LineStyle = 1

$idxRange = $newSheet.Range("A2:A3")
$idxRange.Font.Bold = $true
$idxRange.HorizontalAlignment = -4108  # xlCenter
$idxRange.VerticalAlignment = -4160    # xlTop
$idxRange.Borders.LineStyle = 1
